# Auto-generated Excel COM-interop script
# Applies the 2022-06-13 daily crime data update across all affected worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 2982
$ws.Range('I3').Value = 3063
$ws.Range('G4').Value = 1431
$ws.Range('I4').Value = 730
$ws.Range('I6').Value = 3500
$ws.Range('G7').Value = 24653
$ws.Range('I7').Value = 10551

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I3').Value = 24
$ws.Range('I7').Value = 117

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('I3').Value = 14
$ws.Range('I7').Value = 54

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('I6').Value = 11
$ws.Range('I7').Value = 35

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I3').Value = 103
$ws.Range('G4').Value = 36
$ws.Range('I4').Value = 24
$ws.Range('I6').Value = 92
$ws.Range('G7').Value = 829
$ws.Range('I7').Value = 339

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I2').Value = 52
$ws.Range('I6').Value = 55
$ws.Range('I7').Value = 190

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I2').Value = 97
$ws.Range('I6').Value = 145
$ws.Range('I7').Value = 418

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('I6').Value = 30
$ws.Range('I7').Value = 94

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('I3').Value = 28
$ws.Range('I6').Value = 24
$ws.Range('I7').Value = 90

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I2').Value = 78
$ws.Range('I6').Value = 75
$ws.Range('I7').Value = 234

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I6').Value = 74
$ws.Range('I7').Value = 354
$ws.Range('I8').Value = 665
$ws.Range('I9').Value = 53
$ws.Range('I14').Value = 54
$ws.Range('I17').Value = 12
$ws.Range('I18').Value = 73
$ws.Range('I19').Value = 289
$ws.Range('I20').Value = 271
$ws.Range('I29').Value = 697
$ws.Range('I30').Value = 35
$ws.Range('I31').Value = 94
$ws.Range('I33').Value = 489
$ws.Range('G37').Value = 829
$ws.Range('I37').Value = 339
$ws.Range('I42').Value = 367
$ws.Range('I43').Value = 98
$ws.Range('I44').Value = 82
$ws.Range('I48').Value = 119
$ws.Range('I50').Value = 47
$ws.Range('I51').Value = 95
$ws.Range('I52').Value = 225
$ws.Range('I53').Value = 119
$ws.Range('I54').Value = 237
$ws.Range('I63').Value = 39
$ws.Range('I65').Value = 234
$ws.Range('I66').Value = 25
$ws.Range('I67').Value = 418
$ws.Range('I71').Value = 29
$ws.Range('I76').Value = 165
$ws.Range('I77').Value = 58
$ws.Range('I79').Value = 264
$ws.Range('I83').Value = 207
$ws.Range('I84').Value = 90
$ws.Range('I85').Value = 480
$ws.Range('I86').Value = 62
$ws.Range('I88').Value = 92
$ws.Range('I89').Value = 117
$ws.Range('I90').Value = 126
$ws.Range('I92').Value = 33
$ws.Range('I94').Value = 93
$ws.Range('I97').Value = 83
$ws.Range('I99').Value = 190
$ws.Range('G101').Value = 24653
$ws.Range('I101').Value = 10551

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I6').Value = 38
$ws.Range('I7').Value = 207

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I6').Value = 159
$ws.Range('I7').Value = 489

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I6').Value = 116
$ws.Range('I7').Value = 237

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I3').Value = 245
$ws.Range('I6').Value = 188
$ws.Range('I7').Value = 697

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I3').Value = 79
$ws.Range('I6').Value = 80
$ws.Range('I7').Value = 289

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('I6').Value = 26
$ws.Range('I7').Value = 82

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I3').Value = 23
$ws.Range('I6').Value = 69
$ws.Range('I7').Value = 119

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I2').Value = 37
$ws.Range('I4').Value = 20
$ws.Range('I6').Value = 66
$ws.Range('I7').Value = 165

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I2').Value = 125
$ws.Range('I3').Value = 193
$ws.Range('I7').Value = 480

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('I2').Value = 31
$ws.Range('I6').Value = 13
$ws.Range('I7').Value = 74

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I3').Value = 123
$ws.Range('I7').Value = 367

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I3').Value = 86
$ws.Range('I7').Value = 264

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I2').Value = 75
$ws.Range('I3').Value = 78
$ws.Range('I6').Value = 94
$ws.Range('I7').Value = 271

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('I6').Value = 36
$ws.Range('I7').Value = 73

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('I3').Value = 4
$ws.Range('I7').Value = 12

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I3').Value = 82
$ws.Range('I4').Value = 26
$ws.Range('I6').Value = 51
$ws.Range('I7').Value = 225

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I6').Value = 51
$ws.Range('I7').Value = 93

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('I3').Value = 13
$ws.Range('I7').Value = 47

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('I3').Value = 6
$ws.Range('I7').Value = 25

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('I6').Value = 13
$ws.Range('I7').Value = 53

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('I6').Value = 47
$ws.Range('I7').Value = 83

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('I6').Value = 15
$ws.Range('I7').Value = 33

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I2').Value = 20
$ws.Range('I3').Value = 31
$ws.Range('I7').Value = 92

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 206
$ws.Range('I3').Value = 184
$ws.Range('I4').Value = 43
$ws.Range('I6').Value = 212
$ws.Range('I7').Value = 665

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('I2').Value = 11
$ws.Range('I6').Value = 9
$ws.Range('I7').Value = 62

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('I6').Value = 47
$ws.Range('I7').Value = 126

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I3').Value = 28
$ws.Range('I6').Value = 43
$ws.Range('I7').Value = 95

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('I2').Value = 18
$ws.Range('I3').Value = 16
$ws.Range('I6').Value = 56
$ws.Range('I7').Value = 98

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('I2').Value = 24
$ws.Range('I7').Value = 119

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('I6').Value = 10
$ws.Range('I7').Value = 29

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('I6').Value = 16
$ws.Range('I7').Value = 58

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I2').Value = 120
$ws.Range('I3').Value = 104
$ws.Range('I6').Value = 92
$ws.Range('I7').Value = 354
